# Realestate Update resale numbers 2023-05-31 19:36
# Appends a new data row (row 12) to the CityResaleNum sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Date / Time / Weekday / Week columns are stored as text in this sheet
# (see row 11), so force them to text via a leading apostrophe and then
# reset the cell style back to the same (default) style used by the rest
# of the data rows, so Excel's automatic date/number detection doesn't
# leave a lingering number-format style behind.
$ws.Range("A12").Value = "'2023-05-31"
$ws.Range("A12").Style = $ws.Range("A11").Style

$ws.Range("B12").Value = "19:35:13"
$ws.Range("B12").Style = $ws.Range("B11").Style

$ws.Range("C12").Value = "Wednesday"
$ws.Range("C12").Style = $ws.Range("C11").Style

$ws.Range("D12").Value = "'22"
$ws.Range("D12").Style = $ws.Range("D11").Style

# Remaining columns are plain numbers.
$ws.Range("E12").Value = 120245
$ws.Range("F12").Value = 133429
$ws.Range("G12").Value = 158933
$ws.Range("H12").Value = 131036
$ws.Range("I12").Value = 174858
$ws.Range("J12").Value = 113868
$ws.Range("K12").Value = 198821
$ws.Range("L12").Value = 220315
$ws.Range("M12").Value = 172024
$ws.Range("N12").Value = 120097
$ws.Range("O12").Value = 38710
$ws.Range("P12").Value = 34887
$ws.Range("Q12").Value = 50637
$ws.Range("R12").Value = -1
$ws.Range("S12").Value = 36878
$ws.Range("T12").Value = -1
